# Update "想去人数" (number of people interested) figures for two events
# on both the "展览" sheet and the combined "全部类型" sheet, matching the
# regenerated site data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 423
$wsExpo.Range("F3").Value = 2781
$wsExpo.Range("F4").Value = 126

# --- 全部类型 sheet (aggregated view) ----------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 423
$wsAll.Range("F7").Value = 2781
$wsAll.Range("F8").Value = 126
